$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header values: "Valor Mora" total and "Cant. Periodos" count
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 2160000
$ws.Range("F13").Value = 6

# ---------------------------------------------------------------------------
# 2) Footer block: a new row is inserted before the signature lines, so the
#    old row 25 (signature line) becomes row 26 and old row 26 (legal rep
#    name / signature caption) becomes row 27.
#    Shift bottom-up, copying format and values separately (keeps existing
#    cell styles instead of synthesizing new ones).
# ---------------------------------------------------------------------------
$ws.Range("B26:C26").Copy()
$ws.Range("B27:C27").PasteSpecial(-4122)
$ws.Range("B26:C26").Copy()
$ws.Range("B27:C27").PasteSpecial(-4163)

$ws.Range("H26:J26").Copy()
$ws.Range("H27:J27").PasteSpecial(-4122)
$ws.Range("H26:J26").Copy()
$ws.Range("H27:J27").PasteSpecial(-4163)

$ws.Range("B25:C25").Copy()
$ws.Range("B26:C26").PasteSpecial(-4122)
$ws.Range("B25:C25").Copy()
$ws.Range("B26:C26").PasteSpecial(-4163)

$ws.Range("H25:J25").Copy()
$ws.Range("H26:J26").PasteSpecial(-4122)
$ws.Range("H25:J25").Copy()
$ws.Range("H26:J26").PasteSpecial(-4163)

# Clear the now-stale row 25 (its content moved to row 26)
$ws.Range("B25:C25").ClearContents()
$ws.Range("B25:C25").ClearFormats()
$ws.Range("H25:J25").ClearContents()
$ws.Range("H25:J25").ClearFormats()

# ---------------------------------------------------------------------------
# 3) Data table: add a new last row (21) below the current last row (20),
#    reusing row 20's special "closing border" formatting, then turn row 20
#    into a regular data row (same formatting as rows 16-19).
# ---------------------------------------------------------------------------
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4163)

$ws.Range("B16:J16").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# New period values (most recent first): rows 16..21 -> 2507..2502
$ws.Range("E16").Value = "2507"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2505"
$ws.Range("E19").Value = "2504"
$ws.Range("E20").Value = "2503"
$ws.Range("E21").Value = "2502"

Write-Output "done"
